$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.258580684661865
$ws.Range("B1").Value = 2.375397682189941
$ws.Range("C1").Value = 3.425321340560913
$ws.Range("D1").Value = 2.59455943107605
$ws.Range("E1").Value = 1.355193376541138
